# Updated cryptos list with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row
# with the latest scraped values, preserving their original text format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.668.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.032"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.031"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4386"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3793"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8823"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.857.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.498"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.704"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07154"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.037"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009051"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.030"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.675.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.295"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.081.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.055"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.309"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09049"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7706"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.206"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.005"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.554"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.033"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.149"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01971"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05267"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.843"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5173"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1669"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.870"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06613"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.034"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.699"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4694"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.895"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
